$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 516.0566398604303
$ws.Range("D2").Value = 117.0786055004925
$ws.Range("G2").Value = 472
$ws.Range("C3").Value = 37.3764503358322
$ws.Range("D3").Value = 6.484512362060585
$ws.Range("F3").Value = 32.71
$ws.Range("G3").Value = 37.59
$ws.Range("H3").Value = 41.75
$ws.Range("C4").Value = 2.049962455625566
$ws.Range("D4").Value = 2.527692175702072
$ws.Range("F4").Value = 0.67
$ws.Range("G4").Value = 1.32
$ws.Range("H4").Value = 2.52
$ws.Range("C5").Value = 322.6127146072568
$ws.Range("D5").Value = 8.557186073890973
$ws.Range("F5").Value = 317.69
$ws.Range("G5").Value = 323.15
$ws.Range("H5").Value = 328.38
$ws.Range("C6").Value = 23.74303325768589
$ws.Range("D6").Value = 3.700514144655235
$ws.Range("F6").Value = 21.06
$ws.Range("H6").Value = 26.29
$ws.Range("C7").Value = -75.18715237466306
$ws.Range("D7").Value = 22.09337141334001
$ws.Range("G7").Value = -71
$ws.Range("H7").Value = -57
$ws.Range("C8").Value = 7.954853436707918
$ws.Range("D8").Value = 6.536677735626419
$ws.Range("C9").Value = 9.117618181906218
$ws.Range("D9").Value = 1.603902592242946
$ws.Range("C10").Value = 867.8230302332785
$ws.Range("D10").Value = 0.4610855562334063
$ws.Range("C11").Value = 0.4738713281166023
$ws.Range("D11").Value = 0.5349304607510091
$ws.Range("C12").Value = 22.74796290109255
$ws.Range("D12").Value = 12.29667879762277
$ws.Range("C13").Value = 0.6717072224127997
$ws.Range("D13").Value = 0.7500108620957644
$ws.Range("C14").Value = 1.831115496254315
$ws.Range("D14").Value = 1.669134238883229
$ws.Range("C15").Value = 92.44715237466252
$ws.Range("D15").Value = 22.09337141331377
$ws.Range("F15").Value = 74.25999999999999
$ws.Range("G15").Value = 88.25999999999999
$ws.Range("C16").Value = -84.60145774217064
$ws.Range("D16").Value = 19.94693935633377
$ws.Range("F16").Value = -100.5961208798061
$ws.Range("G16").Value = -82.23249407632485
$ws.Range("H16").Value = -68.59612087980607
$ws.Range("C17").Value = -76.64660430546279
$ws.Range("D17").Value = 24.31944969061291
$ws.Range("F17").Value = -91.34699179957641
$ws.Range("G17").Value = -71.41392685158225
$ws.Range("H17").Value = -57.75746206410165